$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.263052213649928
$ws.Range("C2").Value = 0.3889961846971346
$ws.Range("D2").Value = 0.01665883252591627
$ws.Range("F2").Value = 0.829255421700779
$ws.Range("G2").Value = 0.002415931768689146
$ws.Range("I2").Value = 0.4885738491782838
$ws.Range("B3").Value = 1.12131573768238
$ws.Range("C3").Value = 0.3408341730425946
$ws.Range("D3").Value = 0.01685731912639454
$ws.Range("F3").Value = 0.8198289288857126
$ws.Range("G3").Value = 0.002420730270256
$ws.Range("I3").Value = 0.4965952470778205
$ws.Range("B4").Value = 1.034401361116466
$ws.Range("C4").Value = 0.3112853355500818
$ws.Range("D4").Value = 0.01699604700388946
$ws.Range("F4").Value = 0.8150624508432287
$ws.Range("G4").Value = 0.002423828674492908
$ws.Range("I4").Value = 0.50216040146735
$ws.Range("B5").Value = 0.9990113780824572
$ws.Range("C5").Value = 0.2992492301924585
$ws.Range("D5").Value = 0.01705680592191428
$ws.Range("F5").Value = 0.8133747162603555
$ws.Range("G5").Value = 0.002425129677623175
$ws.Range("I5").Value = 0.5045880485538454
$ws.Range("B6").Value = 0.9931366126592138
$ws.Range("C6").Value = 0.2972509503332788
$ws.Range("D6").Value = 0.01706714979904689
$ws.Range("F6").Value = 0.81310978470907
$ws.Range("G6").Value = 0.002425348030029978
$ws.Range("I6").Value = 0.505000776126959
$ws.Range("B7").Value = 1.033923964587586
$ws.Range("C7").Value = 0.3111229916046057
$ws.Range("D7").Value = 0.01699684932525969
$ws.Range("F7").Value = 0.8150386613430101
$ws.Range("G7").Value = 0.002423846064716095
$ws.Range("I7").Value = 0.5021924960359065
$ws.Range("B8").Value = 1.214158132896159
$ws.Range("C8").Value = 0.3723848916544057
$ws.Range("D8").Value = 0.01672376434665424
$ws.Range("F8").Value = 0.8257918384315275
$ws.Range("G8").Value = 0.002417554801196153
$ws.Range("I8").Value = 0.4912061202129401
$ws.Range("B9").Value = 1.568506119212714
$ws.Range("C9").Value = 0.4927275319830073
$ws.Range("D9").Value = 0.01632261865386297
$ws.Range("F9").Value = 0.8550850067940416
$ws.Range("G9").Value = 0.002406418571106037
$ws.Range("I9").Value = 0.4747905633639178
$ws.Range("B10").Value = 1.829454040565565
$ws.Range("C10").Value = 0.5813171251411404
$ws.Range("D10").Value = 0.01611077276218964
$ws.Range("F10").Value = 0.8817535431648906
$ws.Range("G10").Value = 0.002398960507247283
$ws.Range("I10").Value = 0.465926633928909
$ws.Range("B11").Value = 1.948312367271456
$ws.Range("C11").Value = 0.6216674813978216
$ws.Range("D11").Value = 0.01603262505904723
$ws.Range("F11").Value = 0.895034684887662
$ws.Range("G11").Value = 0.002395722996831776
$ws.Range("I11").Value = 0.4626034244054509
$ws.Range("B12").Value = 1.993343382952332
$ws.Range("C12").Value = 0.6369550885585795
$ws.Range("D12").Value = 0.01600567320014079
$ws.Range("F12").Value = 0.9002316937212242
$ws.Range("G12").Value = 0.002394519216657385
$ws.Range("I12").Value = 0.4614482146225427
$ws.Range("B13").Value = 1.983644169495733
$ws.Range("C13").Value = 0.6336622730001409
$ws.Range("D13").Value = 0.01601135996535419
$ws.Range("F13").Value = 0.8991049262225914
$ws.Range("G13").Value = 0.002394777487090843
$ws.Range("I13").Value = 0.4616923989810928
$ws.Range("B14").Value = 1.952016655325622
$ws.Range("C14").Value = 0.6229250424992188
$ws.Range("D14").Value = 0.01603035466148128
$ws.Range("F14").Value = 0.8954588704177411
$ws.Range("G14").Value = 0.002395623516987684
$ws.Range("I14").Value = 0.4625063087371117
$ws.Range("B15").Value = 1.932646752183985
$ws.Range("C15").Value = 0.6163492112984841
$ws.Range("D15").Value = 0.01604233402938959
$ws.Range("F15").Value = 0.8932474716142735
$ws.Range("G15").Value = 0.002396144621510581
$ws.Range("I15").Value = 0.4630183313940606
$ws.Range("B16").Value = 1.821689432965286
$ws.Range("C16").Value = 0.5786811965588754
$ws.Range("D16").Value = 0.01611624831342695
$ws.Range("F16").Value = 0.8809089277437891
$ws.Range("G16").Value = 0.002399175198619344
$ws.Range("I16").Value = 0.4661581886052915
$ws.Range("B17").Value = 1.753659659957918
$ws.Range("C17").Value = 0.5555864042799499
$ws.Range("D17").Value = 0.01616627377427804
$ws.Range("F17").Value = 0.8736357704374029
$ws.Range("G17").Value = 0.00240107402247848
$ws.Range("I17").Value = 0.4682668845653382
$ws.Range("B18").Value = 1.714545012456028
$ws.Range("C18").Value = 0.542307611066235
$ws.Range("D18").Value = 0.016196760934438
$ws.Range("F18").Value = 0.8695604920745268
$ws.Range("G18").Value = 0.002402180791185455
$ws.Range("I18").Value = 0.4695464130008702
$ws.Range("B19").Value = 1.701303914562004
$ws.Range("C19").Value = 0.5378124344204025
$ws.Range("D19").Value = 0.01620737710518227
$ws.Range("F19").Value = 0.868199156465792
$ws.Range("G19").Value = 0.002402558038039361
$ws.Range("I19").Value = 0.4699910503772742
$ws.Range("B20").Value = 1.760900066965689
$ws.Range("C20").Value = 0.5580443919355389
$ws.Range("D20").Value = 0.01616077095730972
$ws.Range("F20").Value = 0.8743988101665394
$ws.Range("G20").Value = 0.002400870377656718
$ws.Range("I20").Value = 0.4680355018868312
$ws.Range("B21").Value = 1.961305826071907
$ws.Range("C21").Value = 0.6260786112092092
$ws.Range("D21").Value = 0.01602470361329011
$ws.Range("F21").Value = 0.8965252338263099
$ws.Range("G21").Value = 0.002395374416149409
$ws.Range("I21").Value = 0.4622644323559584
$ws.Range("B22").Value = 2.092410590016698
$ws.Range("C22").Value = 0.6705886513340715
$ws.Range("D22").Value = 0.01595117981356609
$ws.Range("F22").Value = 0.9119649201739435
$ws.Range("G22").Value = 0.002391911798716692
$ws.Range("I22").Value = 0.4590949239670721
$ws.Range("B23").Value = 2.022425744332566
$ws.Range("C23").Value = 0.6468284365092245
$ws.Range("D23").Value = 0.01598900424089678
$ws.Range("F23").Value = 0.9036341024482937
$ws.Range("G23").Value = 0.002393748071308913
$ws.Range("I23").Value = 0.4607310296751024
$ws.Range("B24").Value = 1.757626689632389
$ws.Range("C24").Value = 0.5569331397728661
$ws.Range("D24").Value = 0.01616325340349789
$ws.Range("F24").Value = 0.8740535094347024
$ws.Range("G24").Value = 0.002400962398369356
$ws.Range("I24").Value = 0.4681399007436511
$ws.Range("B25").Value = 1.472542850331422
$ws.Range("C25").Value = 0.4601445790254957
$ws.Range("D25").Value = 0.01641665885678734
$ws.Range("F25").Value = 0.8462670363020806
$ws.Range("G25").Value = 0.00240930351866838
$ws.Range("I25").Value = 0.4786751266139788
